$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.769784967280688
$ws.Range("C3").Value = 0.7690550063412681
$ws.Range("D3").Value = 0.7219024138625114

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.7675018516749116
$ws.Range("C4").Value = 0.7706067448399195
$ws.Range("D4").Value = 0.7816142035602641

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.7699470540242878
$ws.Range("C5").Value = 0.8460191610583662
$ws.Range("D5").Value = 0.8033925484733374
